# Finish add visit with insurance test
# - Update the ZipCode test-data value on both the Account and Patient
#   sheets from 90210 to 94105.
# - Leave the cursor on the Account sheet (Email cell) instead of the
#   Patient sheet, matching where the author left off editing.

$wb = $excel.ActiveWorkbook

$wsAccount = $wb.Worksheets.Item("Account")
$wsPatient = $wb.Worksheets.Item("Patient")

# ZipCode value update (same field on both Account + Patient sheets).
$wsAccount.Range("B7").Value = 94105
$wsPatient.Range("B7").Value = 94105

# Leave a selection parked on Patient!B8 (HasMedicareMedicaid row) before
# switching back to Account, matching the saved cursor/tab state.
$wsPatient.Activate() | Out-Null
$wsPatient.Range("B8").Select() | Out-Null

$wsAccount.Activate() | Out-Null
$wsAccount.Range("B4").Select() | Out-Null
